$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update financial figures to corrected (restated) values
$ws.Range("D2").Value = 2362
$ws.Range("E2").Value = 80
$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 72
$ws.Range("H2").Value = 72
$ws.Range("I2").Value = 59
$ws.Range("J2").Value = 13
$ws.Range("K2").Value = 1155
$ws.Range("L2").Value = 747
$ws.Range("M2").Value = 408
$ws.Range("N2").Value = 393
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 79
$ws.Range("Q2").Value = 202
$ws.Range("R2").Value = -176
$ws.Range("S2").Value = 153
$ws.Range("T2").Value = 117
$ws.Range("U2").Value = 85
$ws.Range("V2").Value = 202
$ws.Range("W2").Value = 3.38
$ws.Range("X2").Value = 3.05
$ws.Range("Y2").Value = 17.53
$ws.Range("Z2").Value = 7.56
$ws.Range("AA2").Value = 183.29
$ws.Range("AB2").Value = 398.46
$ws.Range("AC2").Value = 94
$ws.Range("AD2").Value = 29.15
$ws.Range("AE2").Value = 623
$ws.Range("AF2").Value = 4.4
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 63135819

# Row 3: update financial figures to corrected (restated) values
$ws.Range("D3").Value = 2889
$ws.Range("E3").Value = 41
$ws.Range("F3").Value = 41
$ws.Range("G3").Value = 27
$ws.Range("H3").Value = 27
$ws.Range("I3").Value = 22
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 1467
$ws.Range("L3").Value = 1033
$ws.Range("M3").Value = 434
$ws.Range("N3").Value = 415
$ws.Range("O3").Value = 19
$ws.Range("P3").Value = 79
$ws.Range("Q3").Value = 327
$ws.Range("R3").Value = -251
$ws.Range("S3").Value = -10
$ws.Range("T3").Value = 123
$ws.Range("U3").Value = 204
$ws.Range("V3").Value = 212
$ws.Range("W3").Value = 1.41
$ws.Range("X3").Value = 0.93
$ws.Range("Y3").Value = 5.54
$ws.Range("Z3").Value = 2.05
$ws.Range("AA3").Value = 238.05
$ws.Range("AB3").Value = 426.15
$ws.Range("AC3").Value = 36
$ws.Range("AD3").Value = 54.34
$ws.Range("AE3").Value = 658
$ws.Range("AF3").Value = 2.93
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 63135819

# Row 4: update financial figures to corrected (restated) values
$ws.Range("D4").Value = 4119
$ws.Range("E4").Value = 142
$ws.Range("F4").Value = 142
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = 33
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1756
$ws.Range("L4").Value = 1264
$ws.Range("M4").Value = 492
$ws.Range("N4").Value = 473
$ws.Range("O4").Value = 19
$ws.Range("P4").Value = 81
$ws.Range("Q4").Value = 399
$ws.Range("R4").Value = -332
$ws.Range("S4").Value = -104
$ws.Range("T4").Value = 89
$ws.Range("U4").Value = 310
$ws.Range("V4").Value = 170
$ws.Range("W4").Value = 3.45
$ws.Range("X4").Value = 0.87
$ws.Range("Y4").Value = 7.41
$ws.Range("Z4").Value = 2.23
$ws.Range("AA4").Value = 257.15
$ws.Range("AB4").Value = 465.62
$ws.Range("AC4").Value = 51
$ws.Range("AD4").Value = 33.59
$ws.Range("AE4").Value = 727
$ws.Range("AF4").Value = 2.36
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 65001731

# Row 5: update financial figures to corrected (restated) values
$ws.Range("D5").Value = 6128
$ws.Range("E5").Value = 455
$ws.Range("F5").Value = 455
$ws.Range("G5").Value = 421
$ws.Range("H5").Value = 371
$ws.Range("I5").Value = 297
$ws.Range("J5").Value = 74
$ws.Range("K5").Value = 3014
$ws.Range("L5").Value = 2089
$ws.Range("M5").Value = 925
$ws.Range("N5").Value = 835
$ws.Range("O5").Value = 90
$ws.Range("P5").Value = 355
$ws.Range("Q5").Value = 1081
$ws.Range("R5").Value = -162
$ws.Range("S5").Value = -4
$ws.Range("T5").Value = 136
$ws.Range("U5").Value = 946
$ws.Range("V5").Value = 88
$ws.Range("W5").Value = 7.43
$ws.Range("X5").Value = 6.06
$ws.Range("Y5").Value = 45.49
$ws.Range("Z5").Value = 15.57
$ws.Range("AA5").Value = 225.88
$ws.Range("AB5").Value = 135.33
$ws.Range("AC5").Value = 433
$ws.Range("AD5").Value = 8.06
$ws.Range("AE5").Value = 1177
$ws.Range("AF5").Value = 2.97
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 70959464

# Row 6: update financial figures to corrected (restated) values
$ws.Range("D6").Value = 7375
$ws.Range("E6").Value = 433
$ws.Range("F6").Value = 433
$ws.Range("G6").Value = 428
$ws.Range("H6").Value = 338
$ws.Range("I6").Value = 248
$ws.Range("K6").Value = 4975
$ws.Range("L6").Value = 2126
$ws.Range("M6").Value = 2849
$ws.Range("N6").Value = 1852
$ws.Range("P6").Value = 377
$ws.Range("Q6").Value = 275
$ws.Range("R6").Value = -762
$ws.Range("S6").Value = 1518
$ws.Range("T6").Value = 263
$ws.Range("U6").Value = 12
$ws.Range("V6").Value = 21
$ws.Range("W6").Value = 5.87
$ws.Range("X6").Value = 4.59
$ws.Range("Y6").Value = 18.49
$ws.Range("Z6").Value = 8.470000000000001
$ws.Range("AA6").Value = 74.63
$ws.Range("AB6").Value = 390.99
$ws.Range("AC6").Value = 344
$ws.Range("AD6").Value = 6.93
$ws.Range("AE6").Value = 2455
$ws.Range("AF6").Value = 0.97
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 75437074

# Row 6: AI6 no longer reported -> clear the cell entirely
$ws.Range("AI6").ClearContents()

# Rows 7-9: these forecast years are removed from the dataset entirely,
# clear all data cells but keep the row label columns (A, B, C)
$ws.Range("D7:AJ9").ClearContents()
